# Generate Report for Handback
# Updates the localization-status report to reflect a failed handback
# transform for the a557c567-... file in both the zh-cn and de-de
# language sheets (and rolls that status up into the Overview sheet).

$wb = $excel.ActiveWorkbook

$statusFailed = "Handback transform failed"

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3:F3").Value = $statusFailed

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusFailed
$wsZhCn.Range("P3").Value = "Handback file name: zchte225.z5d is different with handoff file name: a557c567-6ba8-4d96-82ea-f127465891bb.7f1a716f99f3b3bc0ddaacbcec273d8a56ce05af.zh-cn."
$wsZhCn.Range("P1").EntireColumn.ColumnWidth = 39.17

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusFailed
$wsDeDe.Range("P3").Value = "Handback file name: zchte225.z5d is different with handoff file name: a557c567-6ba8-4d96-82ea-f127465891bb.7f1a716f99f3b3bc0ddaacbcec273d8a56ce05af.de-de."
$wsDeDe.Range("P1").EntireColumn.ColumnWidth = 39.17
